$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-6) holds the "Förändrad" date; advance it by one day
# (45174 -> 45175, i.e. 2023-09-05 -> 2023-09-06) for each record.
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45175
}
